$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "dylankato" row entirely (row 2), shifting the rows below it up.
$ws.Rows.Item(2).Delete()

# Reorder the remaining two usernames so "leomessi" comes before "sujal.incognito".
$ws.Range("A3").Value = "leomessi"
$ws.Range("A4").Value = "sujal.incognito"

# Update the active selection to match the post-edit state.
$ws.Range("A2").Select()
